# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell (AC1)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Find the last used row of data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD
    $ws.Cells.Item($r, 31).Value = 69   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
